$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the date and volume values between row 2 and row 5
$ws.Range("D2").Value = 44277
$ws.Range("J2").Value = 150

$ws.Range("D5").Value = 44291
$ws.Range("J5").Value = 30
